$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: mark the f12bf0be file as handed back ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = $newStatus
$ws.Range("C3").Value = $newStatus

# --- zh-cn sheet: update status + handback datetime ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B3").Value = $newStatus
$ws.Range("G2").Value = "2016-02-16 10:39:12"
$ws.Range("G3").Value = "2016-02-16 10:39:12"

# --- de-de sheet: update status + handback datetime ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B3").Value = $newStatus
$ws.Range("G2").Value = "2016-02-16 10:39:39"
$ws.Range("G3").Value = "2016-02-16 10:39:39"
